# Update absenteeism data rows 2-11 with new synchronized values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=27319; B="Maria Vitória Pinto";      C="Recursos Humanos"; D="Outros";              E=3; F=45099; G=3245.23}
    @{Row=3;  A=9040;  B="Sarah da Paz";              C="TI";               D="Viagem de negócios";  E=5; F=45088; G=12100.9}
    @{Row=4;  A=31123; B="Isaac da Paz";              C="Vendas";           D="Doença";               E=4; F=45083; G=8106.84}
    @{Row=5;  A=51424; B="Marcelo Peixoto";           C="Financeiro";       D="Problemas pessoais";  E=6; F=45084; G=8298.67}
    @{Row=6;  A=5807;  B="Mariana Lima";              C="Marketing";        D="Outros";              E=7; F=45094; G=5700.06}
    @{Row=7;  A=28831; B="Srta. Luana Nascimento";    C="Vendas";           D="Consulta médica";     E=4; F=45085; G=8820.22}
    @{Row=8;  A=70367; B="Maria Luiza Oliveira";      C="Vendas";           D="Outros";              E=8; F=45081; G=12435.89}
    @{Row=9;  A=25419; B="Letícia Costa";             C="Engenharia";       D="Problemas pessoais";  E=5; F=45095; G=2905.08}
    @{Row=10; A=98709; B="Heitor Costela";            C="Jurídico";         D="Outros";              E=5; F=45102; G=6606.79}
    @{Row=11; A=79812; B="Helena Pires";              C="Recursos Humanos"; D="Doença";               E=7; F=45081; G=2823.48}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
